$wb = $excel.ActiveWorkbook

# Rows whose "latest handoff/handback" timestamp gets refreshed by this
# report-generation run: row 7 ("Handback transform failed") and rows
# 10-16 ("Ready for handoff").
$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

# --- Overview sheet: column D = "Latest Handoff Date" ---
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Range("D$r").Value = "2016-23-20 04:23:34"
}

# --- zh-cn sheet: column E = "Latest Handoff Datetime" ---
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Range("E$r").Value = "2016-03-20 04:23:30"
}

# --- de-de sheet: column E = "Latest Handoff Datetime" ---
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Range("E$r").Value = "2016-03-20 04:23:34"
}
